# TimeLog.xlsx update: add a new day entry (row 18) to the time log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings (added in this edit, in this exact order) ---
# 29: "Figured out dictionary error. Began work on the AI bot and player choices"
# 30: "Figured out dictionary error. Work on the AI bot and player choices."
# 31: "AI can now know when it hits a ship. It can find the neighbours of that tile."
# 32: "AI can now know when it hits a ship. It can find the neighbours of that tile. Need to make sure
#      that AI knows which its has already chosen so it doesnt repeat choices."
#
# Only strings 30 and 32 end up referenced by a cell (D18/E18); 29 and 31 are earlier drafts of the
# same notes that were typed and then edited, which is why they remain in the shared string table.

$ws.Range("Z1").Value = "Figured out dictionary error. Began work on the AI bot and player choices"
$ws.Range("Z2").Value = "Figured out dictionary error. Work on the AI bot and player choices."
$ws.Range("Z3").Value = "AI can now know when it hits a ship. It can find the neighbours of that tile."
$ws.Range("Z4").Value = "AI can now know when it hits a ship. It can find the neighbours of that tile. Need to make sure that AI knows which its has already chosen so it doesnt repeat choices."
$ws.Range("Z1:Z4").ClearContents()

# --- New row 18: new day logged ---
$ws.Range("B18").Value = 45372
$ws.Range("B18").NumberFormat = "d-mmm"
$ws.Range("B18").Font.Name = "Calibri"
$ws.Range("B18").Font.Size = 11
$ws.Range("B18").Borders(7).LineStyle = 1
$ws.Range("B18").Borders(7).LineStyle = -4142
$ws.Range("B18").VerticalAlignment = -4107

$ws.Range("C18").Value = 3

$ws.Range("D18").Value = "Figured out dictionary error. Work on the AI bot and player choices."
$ws.Range("E18").Value = "AI can now know when it hits a ship. It can find the neighbours of that tile. Need to make sure that AI knows which its has already chosen so it doesnt repeat choices."

$ws.Rows("18").RowHeight = 28

# --- Recalculate totals now that a new entry exists ---
$ws.Range("C3").Formula = "=SUM(C7:C39)"
$ws.Range("C4").Formula = "=C2-C3"

# --- Move the active selection down, like a user who just finished typing row 18 ---
$ws.Range("D19").Select()

$wb.Save()
